$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Row 2 (ANI): the "articleBody" class value becomes "content count-br"
$ws.Range("F2").Value = "content count-br"

# 2. Row 3 (TOI): the "articleBody" class value becomes a brand-new unique value "_s30J clearfix"
$ws.Range("F3").Value = "_s30J clearfix"

# 3. Row 7 (TV9) picks up an explicit black Calibri font (instead of the themed default)
#    to match the rest of the data rows, and its row height is normalized to 18.75
#    (same as every other data row) instead of the odd 17.25.
$row7 = $ws.Range("A7:G7")
$row7.Font.Color = 0

$ws.Rows.Item(7).RowHeight = 18.75
